# Auto-generated Excel COM-interop script applying Hyperion_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 17 (ALC) - hunk 0
$ws_ALC.Range("H17").Value = 4289.484
$ws_ALC.Range("J17").Value = 4506
$ws_ALC.Range("L17").Value = 13518
$ws_ALC.Range("N17").Value = -13854

# Row 19 (ALC) - hunk 1
$ws_ALC.Range("H19").Value = 3237
$ws_ALC.Range("I19").Value = 1086.8889
$ws_ALC.Range("J19").Value = 5172.1
$ws_ALC.Range("K19").Value = 1086.8889
$ws_ALC.Range("L19").Value = 5172.1
$ws_ALC.Range("M19").Value = -911.8888999999999
$ws_ALC.Range("N19").Value = -5522.1

# Row 51 (ALC) - hunk 2
$ws_ALC.Range("H51").Value = 6895.8184
$ws_ALC.Range("I51").Value = 1400
$ws_ALC.Range("J51").Value = 7157.524
$ws_ALC.Range("K51").Value = 1400
$ws_ALC.Range("L51").Value = 7157.524
$ws_ALC.Range("M51").Value = -916
$ws_ALC.Range("N51").Value = -8125.524

# Row 76 (ALC) - hunk 3
$ws_ALC.Range("H76").Value = 4646.8335
$ws_ALC.Range("I76").Value = 4646.8335
$ws_ALC.Range("K76").Value = 4646.8335
$ws_ALC.Range("M76").Value = -4331.8335

# Row 79 (ALC) - hunk 4
$ws_ALC.Range("H79").Value = 4646.8335
$ws_ALC.Range("I79").Value = 4646.8335
$ws_ALC.Range("K79").Value = 4646.8335
$ws_ALC.Range("M79").Value = -3554.8335

# Row 112 (ALC) - hunk 5
$ws_ALC.Range("H112").Value = 8457.857
$ws_ALC.Range("I112").Value = 850
$ws_ALC.Range("J112").Value = 9725.833000000001
$ws_ALC.Range("K112").Value = 2550
$ws_ALC.Range("L112").Value = 29177.499
$ws_ALC.Range("M112").Value = -1442
$ws_ALC.Range("N112").Value = -31393.499

# Row 113 (ALC) - hunk 6
$ws_ALC.Range("H113").Value = 4400.7896
$ws_ALC.Range("J113").Value = 6229.8184
$ws_ALC.Range("L113").Value = 6229.8184
$ws_ALC.Range("N113").Value = -12737.8184

# Row 125 (ALC) - hunk 7
$ws_ALC.Range("I125").Value = 1777.8572
$ws_ALC.Range("J125").Value = 7939078.5
$ws_ALC.Range("K125").Value = 16000.7148
$ws_ALC.Range("L125").Value = 71451706.5
$ws_ALC.Range("M125").Value = -13540.7148
$ws_ALC.Range("N125").Value = -71456626.5

# Row 136 (ALC) - hunk 8
$ws_ALC.Range("H136").Value = 134497.25
$ws_ALC.Range("J136").Value = 198999.5
$ws_ALC.Range("L136").Value = 198999.5
$ws_ALC.Range("N136").Value = -209199.5

# Row 138 (ALC) - hunk 9
$ws_ALC.Range("H138").Value = 2714.9167
$ws_ALC.Range("I138").Value = 1413.0952
$ws_ALC.Range("J138").Value = 3148.8572
$ws_ALC.Range("K138").Value = 4239.2856
$ws_ALC.Range("L138").Value = 9446.571599999999
$ws_ALC.Range("M138").Value = 900.7143999999998
$ws_ALC.Range("N138").Value = -19726.5716

# Row 139 (ALC) - hunk 10
$ws_ALC.Range("H139").Value = 100000
$ws_ALC.Range("J139").Value = 100000
$ws_ALC.Range("L139").Value = 100000
$ws_ALC.Range("N139").Value = -110280

# Row 140 (ALC) - hunk 11
$ws_ALC.Range("H140").Value = 45446.25
$ws_ALC.Range("J140").Value = 43263.332
$ws_ALC.Range("L140").Value = 43263.332
$ws_ALC.Range("N140").Value = -53623.332

# Row 2 (ARM) - hunk 12
$ws_ARM.Range("H2").Value = 2828587.5
$ws_ARM.Range("I2").Value = 4040419.8
$ws_ARM.Range("K2").Value = 4040419.8
$ws_ARM.Range("M2").Value = -4040306.8

# Row 25 (ARM) - hunk 13
$ws_ARM.Range("H25").Value = 1516.6666
$ws_ARM.Range("I25").Value = 2200
$ws_ARM.Range("J25").Value = 150
$ws_ARM.Range("K25").Value = 2200
$ws_ARM.Range("L25").Value = 150
$ws_ARM.Range("M25").Value = -1798
$ws_ARM.Range("N25").Value = -954

# Row 61 (ARM) - hunk 14
$ws_ARM.Range("H61").Value = 2324
$ws_ARM.Range("I61").Value = 1818.4706
$ws_ARM.Range("K61").Value = 1818.4706
$ws_ARM.Range("M61").Value = -1606.4706

# Row 97 (ARM) - hunk 15
$ws_ARM.Range("H97").Value = 1198272.5
$ws_ARM.Range("I97").Value = 1797056.1
$ws_ARM.Range("K97").Value = 1797056.1
$ws_ARM.Range("M97").Value = -1796560.1

# Row 110 (ARM) - hunk 16
$ws_ARM.Range("H110").Value = 1266301.2
$ws_ARM.Range("I110").Value = 1392922.2
$ws_ARM.Range("K110").Value = 1392922.2
$ws_ARM.Range("M110").Value = -1390877.2

# Row 116 (ARM) - hunk 17
$ws_ARM.Range("H116").Value = 2828587.5
$ws_ARM.Range("I116").Value = 4040419.8
$ws_ARM.Range("K116").Value = 4040419.8
$ws_ARM.Range("M116").Value = -4038125.8

# Row 136 (ARM) - hunk 18
$ws_ARM.Range("H136").Value = 2324
$ws_ARM.Range("I136").Value = 1818.4706
$ws_ARM.Range("K136").Value = 5455.4118
$ws_ARM.Range("M136").Value = -2905.4118

# Row 3 (BSM) - hunk 19
$ws_BSM.Range("H3").Value = 2828587.5
$ws_BSM.Range("I3").Value = 4040419.8
$ws_BSM.Range("K3").Value = 4040419.8
$ws_BSM.Range("M3").Value = -4040305.8

# Row 64 (BSM) - hunk 20
$ws_BSM.Range("H64").Value = 1195.75
$ws_BSM.Range("J64").Value = 1161.5
$ws_BSM.Range("L64").Value = 1161.5
$ws_BSM.Range("N64").Value = -1611.5

# Row 67 (BSM) - hunk 21
$ws_BSM.Range("H67").Value = 1195.75
$ws_BSM.Range("J67").Value = 1161.5
$ws_BSM.Range("L67").Value = 1161.5
$ws_BSM.Range("N67").Value = -2721.5

# Row 80 (BSM) - hunk 22
$ws_BSM.Range("H80").Value = 314.25
$ws_BSM.Range("J80").Value = 422
$ws_BSM.Range("L80").Value = 422
$ws_BSM.Range("N80").Value = -2418

# Row 83 (BSM) - hunk 23
$ws_BSM.Range("H83").Value = 314.25
$ws_BSM.Range("J83").Value = 422
$ws_BSM.Range("L83").Value = 2110
$ws_BSM.Range("N83").Value = -12094

# Row 86 (BSM) - hunk 24
$ws_BSM.Range("H86").Value = 5885501.5
$ws_BSM.Range("I86").Value = 5885501.5
$ws_BSM.Range("K86").Value = 5885501.5
$ws_BSM.Range("M86").Value = -5884378.5

# Row 89 (BSM) - hunk 25
$ws_BSM.Range("H89").Value = 5885501.5
$ws_BSM.Range("I89").Value = 5885501.5
$ws_BSM.Range("K89").Value = 29427507.5
$ws_BSM.Range("M89").Value = -29421891.5

# Row 134 (BSM) - hunk 26
$ws_BSM.Range("H134").Value = 3481.1025
$ws_BSM.Range("I134").Value = 1526.4584
$ws_BSM.Range("J134").Value = 6608.533
$ws_BSM.Range("K134").Value = 4579.3752
$ws_BSM.Range("L134").Value = 19825.599
$ws_BSM.Range("M134").Value = -2044.3752
$ws_BSM.Range("N134").Value = -24895.599

# Row 31 (CRP) - hunk 27
$ws_CRP.Range("H31").Value = 23465.023
$ws_CRP.Range("I31").Value = 1217.5769
$ws_CRP.Range("K31").Value = 1217.5769
$ws_CRP.Range("M31").Value = -922.5769

# Row 34 (CRP) - hunk 28
$ws_CRP.Range("H34").Value = 23465.023
$ws_CRP.Range("I34").Value = 1217.5769
$ws_CRP.Range("K34").Value = 1217.5769
$ws_CRP.Range("M34").Value = -1015.5769

# Row 102 (CRP) - hunk 29
$ws_CRP.Range("H102").Value = 49995.668
$ws_CRP.Range("J102").Value = 49995.668
$ws_CRP.Range("L102").Value = 49995.668
$ws_CRP.Range("N102").Value = -54863.668

# Row 109 (CRP) - hunk 30
$ws_CRP.Range("H109").Value = 23064.666
$ws_CRP.Range("J109").Value = 23064.666
$ws_CRP.Range("L109").Value = 23064.666
$ws_CRP.Range("N109").Value = -25144.666

# Row 132 (CRP) - hunk 31
$ws_CRP.Range("H132").Value = 72592.516
$ws_CRP.Range("I132").Value = 52211
$ws_CRP.Range("K132").Value = 156633
$ws_CRP.Range("M132").Value = -154103

# Row 134 (CRP) - hunk 32
$ws_CRP.Range("H134").Value = 35368.43
$ws_CRP.Range("I134").Value = 61854.734
$ws_CRP.Range("J134").Value = 4807.3076
$ws_CRP.Range("K134").Value = 185564.202
$ws_CRP.Range("L134").Value = 14421.9228
$ws_CRP.Range("M134").Value = -183029.202
$ws_CRP.Range("N134").Value = -19491.9228

# Row 132 (CUL) - hunk 33
$ws_CUL.Range("H132").Value = 2472.5881
$ws_CUL.Range("I132").Value = 1981.125
$ws_CUL.Range("J132").Value = 2909.4443
$ws_CUL.Range("K132").Value = 17830.125
$ws_CUL.Range("L132").Value = 26184.9987
$ws_CUL.Range("M132").Value = -15300.125
$ws_CUL.Range("N132").Value = -31244.9987

# Row 80 (GSM) - hunk 34
$ws_GSM.Range("H80").Value = 843568.8
$ws_GSM.Range("I80").Value = 1436506.4
$ws_GSM.Range("K80").Value = 1436506.4
$ws_GSM.Range("M80").Value = -1435508.4

# Row 83 (GSM) - hunk 35
$ws_GSM.Range("H83").Value = 843568.8
$ws_GSM.Range("I83").Value = 1436506.4
$ws_GSM.Range("K83").Value = 7182532
$ws_GSM.Range("M83").Value = -7177540

# Row 97 (GSM) - hunk 36
$ws_GSM.Range("H97").Value = 1253799.5
$ws_GSM.Range("I97").Value = 1488655.6
$ws_GSM.Range("K97").Value = 1488655.6
$ws_GSM.Range("M97").Value = -1488159.6

# Row 107 (GSM) - hunk 37
$ws_GSM.Range("H107").Value = 1539
$ws_GSM.Range("J107").Value = 992.6667
$ws_GSM.Range("L107").Value = 992.6667
$ws_GSM.Range("N107").Value = -4832.6667

# Row 16 (LTW) - hunk 38
$ws_LTW.Range("H16").Value = 1146
$ws_LTW.Range("I16").Value = 991.93335
$ws_LTW.Range("J16").Value = 1916.3334
$ws_LTW.Range("K16").Value = 991.93335
$ws_LTW.Range("L16").Value = 1916.3334
$ws_LTW.Range("M16").Value = -821.93335
$ws_LTW.Range("N16").Value = -2256.3334

# Row 61 (LTW) - hunk 39
$ws_LTW.Range("H61").Value = 4119115.5
$ws_LTW.Range("I61").Value = 6949205
$ws_LTW.Range("J61").Value = 2621.818
$ws_LTW.Range("K61").Value = 6949205
$ws_LTW.Range("L61").Value = 2621.818
$ws_LTW.Range("M61").Value = -6949003
$ws_LTW.Range("N61").Value = -3025.818

# Row 82 (LTW) - hunk 40
$ws_LTW.Range("H82").Value = 7937651
$ws_LTW.Range("J82").Value = 1333.3334
$ws_LTW.Range("L82").Value = 1333.3334
$ws_LTW.Range("N82").Value = -2055.3334

# Row 85 (LTW) - hunk 41
$ws_LTW.Range("H85").Value = 7937651
$ws_LTW.Range("J85").Value = 1333.3334
$ws_LTW.Range("L85").Value = 1333.3334
$ws_LTW.Range("N85").Value = -3829.3334

# Row 113 (LTW) - hunk 42
$ws_LTW.Range("H113").Value = 4119115.5
$ws_LTW.Range("I113").Value = 6949205
$ws_LTW.Range("J113").Value = 2621.818
$ws_LTW.Range("K113").Value = 6949205
$ws_LTW.Range("L113").Value = 2621.818
$ws_LTW.Range("M113").Value = -6947035
$ws_LTW.Range("N113").Value = -6961.818

# Row 132 (LTW) - hunk 43
$ws_LTW.Range("H132").Value = 6236.8286
$ws_LTW.Range("I132").Value = 5938.4443
$ws_LTW.Range("J132").Value = 6552.7646
$ws_LTW.Range("K132").Value = 17815.3329
$ws_LTW.Range("L132").Value = 19658.2938
$ws_LTW.Range("M132").Value = -15285.3329
$ws_LTW.Range("N132").Value = -24718.2938

# Row 136 (LTW) - hunk 44
$ws_LTW.Range("H136").Value = 73169.27
$ws_LTW.Range("I136").Value = 109100.16
$ws_LTW.Range("J136").Value = 4900.6
$ws_LTW.Range("K136").Value = 327300.48
$ws_LTW.Range("L136").Value = 14701.8
$ws_LTW.Range("M136").Value = -324750.48
$ws_LTW.Range("N136").Value = -19801.8

# Row 100 (WVR) - hunk 45
$ws_WVR.Range("H100").Value = 1191.7858
$ws_WVR.Range("J100").Value = 3330
$ws_WVR.Range("L100").Value = 6660
$ws_WVR.Range("N100").Value = -7742

# Row 126 (WVR) - hunk 46
$ws_WVR.Range("H126").Value = 2390.7144
$ws_WVR.Range("J126").Value = 1565.6666
$ws_WVR.Range("L126").Value = 4696.9998
$ws_WVR.Range("N126").Value = -9636.9998

# Row 132 (WVR) - hunk 47
$ws_WVR.Range("H132").Value = 84244210
$ws_WVR.Range("I132").Value = 111114290
$ws_WVR.Range("J132").Value = 3633962
$ws_WVR.Range("K132").Value = 333342870
$ws_WVR.Range("L132").Value = 10901886
$ws_WVR.Range("M132").Value = -333340340
$ws_WVR.Range("N132").Value = -10906946
